$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing last row (41, phone 71277620) had its phone number stored
# as text; normalize it to a plain number to match the rest of the sheet.
$ws.Cells.Item(41, 1).Value = 71277620

# Append the new payment row (row 42) for phone 79174445.
# Phone numbers are kept as text (leading apostrophe forces text type,
# same as the other "phone" values already in the sheet).
$ws.Cells.Item(42, 1).Value = "'79174445"
$ws.Cells.Item(42, 3).Value = "Cash"
$ws.Cells.Item(42, 4).Value = "2025-08-18T17:42:14"
$ws.Cells.Item(42, 5).Value = 60
$ws.Cells.Item(42, 7).Value = 60
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
